# Auto-generated edit script: updates crypto price/volume data per commit
# "Updated cryptos list on Sun Sep 22 22:18:27 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.690.71'
$ws.Range("E2").Value = '  -0.67%  '
$ws.Range("D3").Value = '2.535.37'
$ws.Range("E3").Value = '  -1.12%  '
$ws.Range("D4").Value = '''0.998'
$ws.Range("E4").Value = '  -0.15%  '
$ws.Range("D5").Value = '''580.77'
$ws.Range("E5").Value = '  -0.83%  '
$ws.Range("D6").Value = '''141.97'
$ws.Range("E6").Value = '  -4.07%  '
$ws.Range("D7").Value = '''0.999'
$ws.Range("E7").Value = '  -0.11%  '
$ws.Range("D8").Value = '''0.582'
$ws.Range("E8").Value = '  -2.58%  '
$ws.Range("D9").Value = '''0.104'
$ws.Range("E9").Value = '  -3.89%  '
$ws.Range("D10").Value = '''5.54'
$ws.Range("E10").Value = '  -1.92%  '
$ws.Range("E11").Value = '  -0.17%  '
$ws.Range("D12").Value = '''0.345'
$ws.Range("E12").Value = '  -3.28%  '
$ws.Range("D13").Value = '''26.66'
$ws.Range("E13").Value = '  -2.90%  '
$ws.Range("D14").Value = '2.986.32'
$ws.Range("E14").Value = '  -1.36%  '
$ws.Range("D15").Value = '62.502.86'
$ws.Range("E15").Value = '  -0.87%  '
$ws.Range("D16").Value = '''0.0000143'
$ws.Range("E16").Value = '  -2.91%  '
$ws.Range("D17").Value = '2.545.48'
$ws.Range("E17").Value = '  -0.65%  '
$ws.Range("D18").Value = '''10.93'
$ws.Range("E18").Value = '  -3.76%  '
$ws.Range("D19").Value = '''336.47'
$ws.Range("E19").Value = '  -2.10%  '
$ws.Range("D20").Value = '''4.26'
$ws.Range("E20").Value = '  -3.97%  '
$ws.Range("D21").Value = '''6.52'
$ws.Range("E21").Value = '  -5.17%  '
$ws.Range("E22").Value = '  -0.08%  '
$ws.Range("D23").Value = '''5.70'
$ws.Range("E23").Value = '  +3.02%  '
$ws.Range("D24").Value = '''67.05'
$ws.Range("E24").Value = '  +0.74%  '
$ws.Range("E25").Value = '  +3.24%  '
$ws.Range("D26").Value = '''1.58'
$ws.Range("E26").Value = '  -2.87%  '
$ws.Range("D27").Value = '''0.162'
$ws.Range("E27").Value = '  -4.59%  '
$ws.Range("E28").Value = '  -0.10%  '
$ws.Range("D29").Value = '''7.82'
$ws.Range("E29").Value = '  -5.04%  '
$ws.Range("D30").Value = '''8.07'
$ws.Range("E30").Value = '  -4.34%  '
$ws.Range("D31").Value = '''1.92'
$ws.Range("E31").Value = '  -3.27%  '
$ws.Range("D32").Value = '''464.05'
$ws.Range("E32").Value = '  +0.81%  '
$ws.Range("D33").Value = '0.0₃0781'
$ws.Range("E33").Value = '  -5.62%  '
$ws.Range("D34").Value = '''1.64'
$ws.Range("E34").Value = '  +1.12%  '
$ws.Range("D35").Value = '''175.25'
$ws.Range("E35").Value = '  -0.27%  '
$ws.Range("E36").Value = '  +0.04%  '
$ws.Range("D37").Value = '''0.393'
$ws.Range("E37").Value = '  -3.06%  '
$ws.Range("D38").Value = '''18.60'
$ws.Range("E38").Value = '  -3.03%  '
$ws.Range("D39").Value = '''4.44'
$ws.Range("E39").Value = '  -2.72%  '
$ws.Range("E40").Value = '  +0.03%  '
$ws.Range("D41").Value = '''1.67'
$ws.Range("E41").Value = '  -4.75%  '
$ws.Range("D42").Value = '''39.99'
$ws.Range("E42").Value = '  +1.20%  '
$ws.Range("D43").Value = '''155.50'
$ws.Range("E43").Value = '  +3.00%  '
$ws.Range("D44").Value = '''3.64'
$ws.Range("E44").Value = '  -5.04%  '
$ws.Range("B45").Value = 'Mantle'
$ws.Range("C45").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D45").Value = '''0.623'
$ws.Range("E45").Value = '  +1.69%  '
$ws.Range("B46").Value = 'InjectiveProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D46").Value = '''20.79'
$ws.Range("E46").Value = '  -0.93%  '
$ws.Range("D47").Value = '''0.0527'
$ws.Range("E47").Value = '  -4.01%  '
$ws.Range("D48").Value = '''0.0951'
$ws.Range("E48").Value = '  -2.55%  '
$ws.Range("D49").Value = '''0.0233'
$ws.Range("E49").Value = '  -3.40%  '
$ws.Range("D50").Value = '''17.74'
$ws.Range("E50").Value = '  -3.93%  '
$ws.Range("D51").Value = '''11.36'
$ws.Range("E51").Value = '  -0.32%  '
